$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new model ("Multinomial Naive Bayes") finished training, so its results
# are appended in column F, mirroring the layout already used for the other
# models in columns B-E. Values are entered in the same logical order the
# training script would have produced them (name, metrics, run id, then the
# per-sample predictions) so new shared-string entries land in that order.
# ---------------------------------------------------------------------------

function Set-TextValue($cell, $text, $formatSourceCell) {
    $cell.Value = "'" + $text
    $formatSourceCell.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# Model name.
$ws.Range("F1").Value = "Multinomial Naive Bayes"

# Cross validation / accuracy / F1 / ROC-AUC summary rows. These are stored
# as text (not numbers) just like the existing B24:E27 block, so a leading
# apostrophe is used to force text entry; the formatting is then restored
# from the neighbouring cell so no stray "quote prefix" style is introduced.
Set-TextValue $ws.Range("F24") "0.8992" $ws.Range("E24")
Set-TextValue $ws.Range("F25") "0.8972" $ws.Range("E25")
Set-TextValue $ws.Range("F26") "0.9018" $ws.Range("E26")
Set-TextValue $ws.Range("F27") "0.9708" $ws.Range("E27")

# Training run id.
$ws.Range("F2").Value = "anxiety_model_20250510_2039"

# Sample texts (rows 3,5,7,...,21) are identical across every model column,
# and the prediction results (rows 4,6,8,...,22) are specific to this model.
$ws.Range("F3").Value  = $ws.Range("E3").Value2
$ws.Range("F4").Value  = [char]0x2192 + "Anksiyete (98.64%)"

$ws.Range("F5").Value  = $ws.Range("E5").Value2
$ws.Range("F6").Value  = [char]0x2192 + "Anksiyete (59.12%)"

$ws.Range("F7").Value  = $ws.Range("E7").Value2
$ws.Range("F8").Value  = [char]0x2192 + "Anksiyete (98.38%)"

$ws.Range("F9").Value  = $ws.Range("E9").Value2
$ws.Range("F10").Value = [char]0x2192 + "Normal (18.99%)"

$ws.Range("F11").Value = $ws.Range("E11").Value2
$ws.Range("F12").Value = [char]0x2192 + "Anksiyete (87.41%)"

$ws.Range("F13").Value = $ws.Range("E13").Value2
$ws.Range("F14").Value = [char]0x2192 + "Anksiyete (83.68%)"

$ws.Range("F15").Value = $ws.Range("E15").Value2
$ws.Range("F16").Value = [char]0x2192 + "Anksiyete (94.38%)"

$ws.Range("F17").Value = $ws.Range("E17").Value2
$ws.Range("F18").Value = [char]0x2192 + "Anksiyete (64.69%)"

$ws.Range("F19").Value = $ws.Range("E19").Value2
$ws.Range("F20").Value = [char]0x2192 + "Anksiyete (96.03%)"

$ws.Range("F21").Value = $ws.Range("E21").Value2
$ws.Range("F22").Value = [char]0x2192 + "Normal (47.84%)"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Font color of each new prediction cell should reflect the verdict, exactly
# like the conditional styling already used in columns B-E (style tied to
# font 3 = red/"Anksiyete", font 4 = green/"Normal").
# ---------------------------------------------------------------------------

$ws.Range("B4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F16").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# The window had been scrolled/zoomed down to review the newly completed
# training block at the bottom of the sheet.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.Zoom = 64
$ws.Range("D24").Select() | Out-Null
